$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.168.73"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.753.07"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5276"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06209"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "1.746.70"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07192"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6495"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.639"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9993"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "26.036.89"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006758"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "1.969.41"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.345"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.771"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.263"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.520"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.822"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08310"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.830"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.646"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.018"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6400"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.53%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01623"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7483"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.040"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1158"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.411"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05351"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3496"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.622"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
